$d = $word.ActiveDocument

# Word Find.Execute constant for "replace all"
$wdReplaceAll = 2
$wdFindContinue = 1

# 1. Update the title heading and the bolded repeat of it later in the doc
#    (two occurrences of the exact same string; replace both in one pass)
$d.Content.Find.Execute(
    "Play Johnan Legendarian for Free - Features Super Bonus and High RTP",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Play Johnan Legendarian for Free - Review and Features", $wdReplaceAll) | Out-Null

# 2. Insert a brand-new bullet ("High volatility and 10 paylines") right
#    before the existing "Super Bonus feature with 10 free spins" bullet.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Super Bonus feature with 10 free spins*") {
        $rng = $p.Range
        $rng.Collapse(1)
        $rng.InsertBefore("High volatility and 10 paylines`r")
        break
    }
}

# 3. "Possibility of winning up to 5000 times your bet" -> "Epic soundtrack and immersive graphics"
$d.Content.Find.Execute(
    "Possibility of winning up to 5000 times your bet",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Epic soundtrack and immersive graphics", $wdReplaceAll) | Out-Null

# 4. "Great graphics and epic soundtrack" -> "Good chances of winning with extra features"
$d.Content.Find.Execute(
    "Great graphics and epic soundtrack",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Good chances of winning with extra features", $wdReplaceAll) | Out-Null

# 5. Remove the whole "High RTP of 96.03%" bullet paragraph entirely
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*High RTP of 96.03%*") {
        $p.Range.Delete()
        break
    }
}

# 6. "Base game winning potential is not as high" -> "Not as high winning potential in the base game"
$d.Content.Find.Execute(
    "Base game winning potential is not as high",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Not as high winning potential in the base game", $wdReplaceAll) | Out-Null

# 7. "Requires patience to unlock bonus features" -> "Requires patience to unlock Super Bonus feature"
$d.Content.Find.Execute(
    "Requires patience to unlock bonus features",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Requires patience to unlock Super Bonus feature", $wdReplaceAll) | Out-Null

# 8. Update the meta-description-style italic paragraph at the end
$d.Content.Find.Execute(
    "Read our review of Johnan Legendarian, a high volatility slot game with 10 paylines, featuring a Super Bonus and possibility of winning up to 5000x your bet. Play for free.",
    $true, $false, $false, $false, $false, $true, $wdFindContinue, $false,
    "Read our review of Johnan Legendarian and play this game for free. Discover its features and winning potential.", $wdReplaceAll) | Out-Null
